# Apply the updated crypto market snapshot (prices + 1h volume deltas) to Sheet1.
# Source data comes straight from the commit's row-level OOXML diff; each entry below
# is written as literal text so values such as "59.472.65" (thousand-dot formatted)
# or "518.11" keep matching the sheet's existing inline-string cell layout instead of
# being auto-parsed into numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '59.472.65' }
    @{ Cell = 'E2'; Value = '  -1.23%  ' }
    @{ Cell = 'D3'; Value = '2.639.70' }
    @{ Cell = 'E3'; Value = '  +0.58%  ' }
    @{ Cell = 'E4'; Value = '  +0.06%  ' }
    @{ Cell = 'D5'; Value = '518.11' }
    @{ Cell = 'E5'; Value = '  -0.84%  ' }
    @{ Cell = 'D6'; Value = '147.26' }
    @{ Cell = 'E6'; Value = '  -1.17%  ' }
    @{ Cell = 'D7'; Value = '0.995' }
    @{ Cell = 'E7'; Value = '  -0.34%  ' }
    @{ Cell = 'D8'; Value = '0.574' }
    @{ Cell = 'E8'; Value = '  +0.46%  ' }
    @{ Cell = 'D9'; Value = '2.663.97' }
    @{ Cell = 'E9'; Value = '  +1.23%  ' }
    @{ Cell = 'D10'; Value = '6.46' }
    @{ Cell = 'E10'; Value = '  +2.26%  ' }
    @{ Cell = 'E11'; Value = '  +1.03%  ' }
    @{ Cell = 'D12'; Value = '0.341' }
    @{ Cell = 'E12'; Value = '  -0.42%  ' }
    @{ Cell = 'E13'; Value = '  -1.55%  ' }
    @{ Cell = 'D14'; Value = '3.109.19' }
    @{ Cell = 'E14'; Value = '  +0.84%  ' }
    @{ Cell = 'D15'; Value = '59.398.60' }
    @{ Cell = 'E15'; Value = '  -1.34%  ' }
    @{ Cell = 'D16'; Value = '21.24' }
    @{ Cell = 'E16'; Value = '  +0.08%  ' }
    @{ Cell = 'E17'; Value = '  +0.29%  ' }
    @{ Cell = 'D18'; Value = '2.666.01' }
    @{ Cell = 'E18'; Value = '  +1.61%  ' }
    @{ Cell = 'D19'; Value = '4.62' }
    @{ Cell = 'E19'; Value = '  -0.43%  ' }
    @{ Cell = 'D20'; Value = '346.78' }
    @{ Cell = 'E20'; Value = '  +1.64%  ' }
    @{ Cell = 'E21'; Value = '  +0.78%  ' }
    @{ Cell = 'D22'; Value = '6.19' }
    @{ Cell = 'E22'; Value = '  +1.03%  ' }
    @{ Cell = 'D23'; Value = '0.998' }
    @{ Cell = 'E23'; Value = '  +0.28%  ' }
    @{ Cell = 'D24'; Value = '61.70' }
    @{ Cell = 'E24'; Value = '  +1.72%  ' }
    @{ Cell = 'D25'; Value = '0.425' }
    @{ Cell = 'E25'; Value = '  +1.11%  ' }
    @{ Cell = 'D26'; Value = '2.777.49' }
    @{ Cell = 'E26'; Value = '  +1.48%  ' }
    @{ Cell = 'D27'; Value = '0.994' }
    @{ Cell = 'E27'; Value = '  -0.03%  ' }
    @{ Cell = 'E28'; Value = '  -0.22%  ' }
    @{ Cell = 'D29'; Value = '0.0₃0821' }
    @{ Cell = 'E29'; Value = '  +1.52%  ' }
    @{ Cell = 'E30'; Value = '  +1.99%  ' }
    @{ Cell = 'D31'; Value = '0.997' }
    @{ Cell = 'E31'; Value = '  -0.33%  ' }
    @{ Cell = 'D32'; Value = '6.51' }
    @{ Cell = 'E32'; Value = '  +8.48%  ' }
    @{ Cell = 'E33'; Value = '  +0.63%  ' }
    @{ Cell = 'E34'; Value = '  -0.15%  ' }
    @{ Cell = 'D35'; Value = '150.27' }
    @{ Cell = 'E35'; Value = '  +0.25%  ' }
    @{ Cell = 'D36'; Value = '1.04' }
    @{ Cell = 'E36'; Value = '  +13.57%  ' }
    @{ Cell = 'E37'; Value = '  +2.75%  ' }
    @{ Cell = 'E38'; Value = '  +3.05%  ' }
    @{ Cell = 'D39'; Value = '0.871' }
    @{ Cell = 'E39'; Value = '  +0.73%  ' }
    @{ Cell = 'D40'; Value = '36.71' }
    @{ Cell = 'E40'; Value = '  +0.58%  ' }
    @{ Cell = 'D41'; Value = '3.74' }
    @{ Cell = 'E41'; Value = '  +3.03%  ' }
    @{ Cell = 'E42'; Value = '  -0.57%  ' }
    @{ Cell = 'D43'; Value = '287.98' }
    @{ Cell = 'E43'; Value = '  -0.42%  ' }
    @{ Cell = 'D44'; Value = '0.618' }
    @{ Cell = 'E44'; Value = '  -1.38%  ' }
    @{ Cell = 'D45'; Value = '0.0994' }
    @{ Cell = 'E45'; Value = '  -1.04%  ' }
    @{ Cell = 'D46'; Value = '0.993' }
    @{ Cell = 'E46'; Value = '  -0.51%  ' }
    @{ Cell = 'D47'; Value = '19.70' }
    @{ Cell = 'E47'; Value = '  +0.98%  ' }
    @{ Cell = 'E48'; Value = '  -0.37%  ' }
    @{ Cell = 'E49'; Value = '  +0.70%  ' }
    @{ Cell = 'B50'; Value = 'RenderToken' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D50'; Value = '4.76' }
    @{ Cell = 'E51'; Value = '  -1.28%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $text = $u.Value
    # Force a leading quote-prefix for strings that look like plain numbers so Excel
    # stores them as text (matching the workbook's existing inline-string cells)
    # instead of silently converting them to a numeric cell value.
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}
